# This workbook is a weekly price log. A new weekly observation was added
# as a new row right after the existing header/first data row (i.e. it
# becomes the new row 3), pushing all the previous observations (formerly
# rows 3-80) down by one row (to rows 4-81).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 3; this shifts rows 3:80 down to 4:81
# and keeps the rest of the sheet (including formatting) intact.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly observation.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 45043
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 38
$ws.Range("K3").Value = 29000
$ws.Range("L3").Value = 29000
$ws.Range("M3").Value = 29000
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 1160
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
